$d = $word.ActiveDocument

# Paragraph 1 = Title "TBV Tags"
# Paragraph 2 = the blank paragraph that must be kept
# Paragraph 3 = first paragraph to remove ("[PUMP:TBV:1111]")
# Paragraph 16 = last paragraph in the body (final bolus test details)
$startPar = $d.Paragraphs(3)
$endPar = $d.Paragraphs($d.Paragraphs.Count)

$rng = $d.Range($startPar.Range.Start, $endPar.Range.End)
$rng.Delete()
